$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "glossary" sheet: update definitions/field types/text formats for the
#    panel_deployment_* and sample_retrieval_* rows.
# ---------------------------------------------------------------------------
$glossary = $wb.Worksheets.Item("glossary")

# Row 11: panel_deployment_day
$glossary.Range("B11").Value = "The day the panel was deployed in the field"
$glossary.Range("C11").Value = "numeric"
$glossary.Range("D11").Value = "DD"
$glossary.Range("F11").Value = "sample metadata"

# Row 12: panel_deployment_month
$glossary.Range("B12").Value = "The month the panel was deployed in the field"
$glossary.Range("C12").Value = "numeric"
$glossary.Range("D12").Value = "MM"
$glossary.Range("F12").Value = "sample metadata"

# Row 13: panel_deployment_year
$glossary.Range("C13").Value = "numeric"
$glossary.Range("D13").Value = "YYYY"
$glossary.Range("F13").Value = "sample metadata"

# Row 22: sample_retrieval_day
$glossary.Range("C22").Value = "numeric"
$glossary.Range("D22").Value = "DD"
$glossary.Range("F22").Value = "sample metadata, biomass data, sessile species data, mobile fauna data, percent cover"

# Row 23: sample_retrieval_month
$glossary.Range("C23").Value = "numeric"
$glossary.Range("D23").Value = "MM"
$glossary.Range("F23").Value = "sample metadata, biomass data, sessile species data, mobile fauna data, percent cover"

# Row 24: sample_retrieval_year
$glossary.Range("C24").Value = "numeric"
$glossary.Range("D24").Value = "YYYY"
$glossary.Range("F24").Value = "sample metadata, biomass data, sessile species data, mobile fauna data, percent cover"

# ---------------------------------------------------------------------------
# 2) "biomass data" sheet: the panel_deployment_* date columns are replaced
#    with the sample_retrieval_* date columns, matching the other data
#    entry sheets.
# ---------------------------------------------------------------------------
$biomass = $wb.Worksheets.Item("biomass data")
$biomass.Range("D1").Value = "sample_retrieval_year"
$biomass.Range("E1").Value = "sample_retrieval_month"
$biomass.Range("F1").Value = "sample_retrieval_day"

# ---------------------------------------------------------------------------
# 3) Add a thin black bottom border under the header row of every data entry
#    sheet (sample metadata, biomass data, sessile species data,
#    mobile fauna data, percent cover).
# ---------------------------------------------------------------------------
$dataSheetNames = @("sample metadata", "biomass data", "sessile species data", "mobile fauna data", "percent cover")
foreach ($name in $dataSheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count
    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
    $headerRange.Borders.Item(9).Color = 0
    $headerRange.Borders.Item(9).Weight = 2
    $headerRange.Borders.Item(9).LineStyle = 1
}
